$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.860646486282349
$ws.Range("B1").Value = 3.747165679931641
$ws.Range("C1").Value = 1.65989089012146
$ws.Range("D1").Value = 1.06648588180542
$ws.Range("E1").Value = 1.101778864860535
